# Fix DVI timing constraint error
#
# Duplicate Sheet1 into a new "Sheet1 (2)" tab (placed right after Sheet1)
# and adjust its horizontal/vertical timing inputs for a DVI-compatible
# pixel clock (25 MHz / 6) instead of the VGA 25 MHz clock used on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Copy Sheet1 and place the copy immediately after it; Excel names it
# "Sheet1 (2)" automatically.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# Pixel frequency becomes a derived value (25 MHz divided by 6) on the new
# sheet instead of a hard-coded 25,000,000.
$ws2.Range("B5").Formula = "=25000000/6"

# Horizontal timing inputs (pixels)
$ws2.Range("B8").Value = 160
$ws2.Range("B9").Value = 80
$ws2.Range("B10").Value = 12
$ws2.Range("B11").Value = 204

# Vertical timing inputs (lines)
$ws2.Range("B15").Value = 144
$ws2.Range("B16").Value = 0
$ws2.Range("B17").Value = 0
$ws2.Range("B18").Value = 10

# Restore the on-screen selections to match the authored state: the new
# sheet's cursor sits at F6, while Sheet1 (which stays the active tab)
# ends up with its cursor at F18.
$ws2.Range("F6").Select()
$ws1.Activate()
$ws1.Range("F18").Select()
